# Populate "Example set of parameters" sheet with the final testing values.
# Cells are written in the same chronological order the original author used
# (reconstructed from the shared-string interning order in the target file)
# so that xl/sharedStrings.xml comes out in the same sequence.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Sampling rate (pts/s)"
$ws.Cells.Item(1, 2).Value = 10

$ws.Cells.Item(2, 1).Value = "Starting temperature (°C)"
$ws.Cells.Item(2, 2).Value = 0

$ws.Cells.Item(3, 1).Value = "Final temperature (°C)"
$ws.Cells.Item(3, 2).Value = 180

$ws.Cells.Item(5, 1).Value = "Heating rate (°C/min)"
$ws.Cells.Item(5, 2).Value = 2

$ws.Cells.Item(6, 1).Value = "Modulation amplitude (°c)"
$ws.Cells.Item(6, 2).Value = 0.212

$ws.Cells.Item(4, 1).Value = "Modulation period (s)"
$ws.Cells.Item(4, 2).Value = 40

$ws.Cells.Item(7, 1).Value = "Phase of the modulated heat flow with respect to the temperature modulation (rad)"
$ws.Cells.Item(7, 2).Value = -0.2

$ws.Cells.Item(8, 1).Value = "Degree of smoothing"
$ws.Cells.Item(8, 2).Value = 0.05

$ws.Cells.Item(9, 1).Value = "Number of Gaussians "
$ws.Cells.Item(9, 2).Value = 3

$ws.Cells.Item(10, 1).Value = "Slope of the reversing heat flow before the Tg (J/(g*°C))"
$ws.Cells.Item(10, 2).Value = -0.0001

$ws.Cells.Item(11, 1).Value = "Slope of the reversing heat flow after the Tg (J/(g*°C))"
$ws.Cells.Item(11, 2).Value = -0.00012

$ws.Cells.Item(12, 1).Value = "Starting value of the reversing heat flow (J/g)"
$ws.Cells.Item(12, 2).Value = -0.04

$ws.Cells.Item(15, 1).Value = "Starting value of the total heat capacity (J/g)"
$ws.Cells.Item(15, 2).Value = 1.05

$ws.Cells.Item(16, 1).Value = "THF Tg start, end, and midpoint separated by commas (°C)"

$ws.Cells.Item(17, 1).Value = "RHF start, end, and midpoint separated by commas (°C)"

$ws.Cells.Item(18, 1).Value = "Jump in heat capacity at the Tg (J/(g*°C))"
$ws.Cells.Item(18, 2).Value = 0.3

$ws.Cells.Item(13, 1).Value = "Slope of the total heat capacity before the Tg (J/(g*°C))"
$ws.Cells.Item(13, 2).Value = 0.0008

$ws.Cells.Item(14, 1).Value = "Slope of the total heat capacity after the Tg (J/(g*°C))"
$ws.Cells.Item(14, 2).Value = 0.0009

$ws.Cells.Item(16, 2).Value = "70,90,80"

$ws.Cells.Item(19, 1).Value = "Gaussian 1: Onset (°C), End (°C), Enthalpy (J/g)"

$ws.Cells.Item(20, 1).Value = "Gaussian 2: Onset (°C), End (°C), Enthalpy (J/g)"

$ws.Cells.Item(21, 1).Value = "Gaussian 3: Onset (°C), End (°C), Enthalpy (J/g)"

$ws.Cells.Item(20, 2).Value = "100,120,1"

$ws.Cells.Item(21, 2).Value = "140,160,-1.2"

$ws.Cells.Item(17, 2).Value = "73,93,83"

$ws.Cells.Item(19, 2).Value = "70,80,-0.3"

# Best-fit column widths (matches the author's final bestFit autofit of
# column A/B based on their content: ~77.57 chars for the longest label in
# column A, ~9.86 chars for the widest numeric value in column B).
$ws.Columns.Item(1).ColumnWidth = 77.5703125
$ws.Columns.Item(2).ColumnWidth = 9.85546875

$ws.Range("A13").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 13

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("B20").Select() | Out-Null
